# Apply updated crypto price/volume figures to the 'cryptos' worksheet.
# Values are taken verbatim (including the 2-space padding around the
# percent figures in column E) so the resulting cell text matches exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.204.60"
$ws.Range("E2").Value = "  +0.84%  "
# Row 3
$ws.Range("D3").Value = "2.317.61"
$ws.Range("E3").Value = "  +0.62%  "
# Row 4
$ws.Range("E4").Value = "  +0.12%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.14%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.76%  "
# Row 7
$ws.Range("E7").Value = "  +0.11%  "
# Row 8
$ws.Range("E8").Value = "  +2.75%  "
# Row 9
$ws.Range("D9").Value = "2.313.11"
$ws.Range("E9").Value = "  +0.49%  "
# Row 10
$ws.Range("E10").Value = "  -0.86%  "
# Row 11
$ws.Range("E11").Value = "  -0.66%  "
# Row 12
$ws.Range("E12").Value = "  +0.73%  "
# Row 13
$ws.Range("E13").Value = "  +0.08%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.20%  "
# Row 15
$ws.Range("D15").Value = "2.729.80"
$ws.Range("E15").Value = "  +0.87%  "
# Row 16
$ws.Range("D16").Value = "59.077.18"
$ws.Range("E16").Value = "  +0.73%  "
# Row 17
$ws.Range("E17").Value = "  +0.29%  "
# Row 18
$ws.Range("D18").Value = "2.342.24"
$ws.Range("E18").Value = "  +2.09%  "
# Row 19
$ws.Range("E19").Value = "  -0.21%  "
# Row 20
$ws.Range("E20").Value = "  -3.25%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.10%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.04%  "
# Row 23
$ws.Range("E23").Value = "  -0.16%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.87%  "
# Row 25
$ws.Range("E25").Value = "  +2.55%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.11%  "
# Row 27
$ws.Range("E27").Value = "  -1.72%  "
# Row 28
$ws.Range("E28").Value = "  -1.02%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.07%  "
# Row 30
$ws.Range("E30").Value = "  +5.92%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.37%  "
# Row 32
$ws.Range("E32").Value = "  +2.29%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.89"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.94%  "
# Row 34
$ws.Range("E34").Value = "  +1.38%  "
# Row 35
$ws.Range("E35").Value = "  +6.43%  "
# Row 37
$ws.Range("E37").Value = "  +0.86%  "
# Row 38
$ws.Range("E38").Value = "  -0.04%  "
# Row 39
$ws.Range("E39").Value = "  +3.46%  "
# Row 40
$ws.Range("E40").Value = "  +1.07%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "304.54"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.98%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.17"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.55%  "
# Row 43
$ws.Range("E43").Value = "  +0.80%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0959"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.91%  "
# Row 45
$ws.Range("E45").Value = "  -0.77%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.558"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.40%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.34%  "
# Row 48
$ws.Range("E48").Value = "  -1.36%  "
# Row 49
$ws.Range("E49").Value = "  -0.05%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.65"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.11%  "
# Row 51
$ws.Range("E51").Value = "  +2.16%  "
